# Two ICDC "startup" query scripts were updated to resolve a wait-time issue:
# the CasesTab Cypher query (column B, row 2) no longer joins/returns the
# Cohort, and the Cases/Samples/Files query cells were re-saved so the
# workbook's shared-string table settles in the order Cases, Samples, Files
# is reflected via their (unchanged) Sample/File query text being re-applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CasesTab query (B2): drop the trailing
#   OPTIONAL MATCH (co:cohort)...  /  coalesce(co.cohort_description,'') AS `Cohort`
# output column from the RETURN clause.
$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Basset Hound'',''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate'', ''Bladder, Urethra'', ''Bladder, Urethra, Prostate'', ''Urethra, Prostate''] and diag.best_response in [''Not Determined'', ''Partial Response'',''Progressive Disease'',''Stable Disease'']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $casesQuery

# SamplesTab query (B3) and FilesTab query (B4) text is unchanged; re-applied
# so the workbook is re-saved with the Cases query's now-orphaned old text
# dropped from the shared-string table.
$sampleQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Basset Hound'',''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate'', ''Bladder, Urethra'', ''Bladder, Urethra, Prostate'', ''Urethra, Prostate''] and diag.best_response in [''Not Determined'', ''Partial Response'',''Progressive Disease'',''Stable Disease'']

 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'
$ws.Range("B3").Value = $sampleQuery

$fileQuery = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Basset Hound'',''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate'', ''Bladder, Urethra'', ''Bladder, Urethra, Prostate'', ''Urethra, Prostate''] and diag.best_response in [''Not Determined'', ''Partial Response'',''Progressive Disease'',''Stable Disease'']
    
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '''') AS `File Name`, 
        coalesce(f.file_type, '''') AS `File Type`, 
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `File Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'
$ws.Range("B4").Value = $fileQuery

# The shorter CasesTab query wraps to one fewer line, so row 2's height
# shrinks from 348 to 333.5.
$ws.Rows.Item(2).RowHeight = 333.5

# Move the active selection (and top-left scroll position) to B2.
$ws.Range("B2").Select() | Out-Null
